$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-157)
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22).
$ws.Range("C2:C157").Value = 45221
